# Auto-generated edit script applying the Adamantoise_Profits market-data refresh
# (chore: update Sheets via scheduled runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 918.21313
$ws.Range("I15").Value = 918.21313
$ws.Range("K15").Value = 2754.63939
$ws.Range("M15").Value = -2585.63939
$ws.Range("H39").Value = 2029.7
$ws.Range("J39").Value = 3296.1667
$ws.Range("L39").Value = 9888.500100000001
$ws.Range("N39").Value = -10480.5001
$ws.Range("H106").Value = 3334604.2
$ws.Range("I106").Value = 3510068.8
$ws.Range("K106").Value = 3510068.8
$ws.Range("M106").Value = -3509437.8
$ws.Range("H107").Value = 1831.35
$ws.Range("I107").Value = 2279.8667
$ws.Range("K107").Value = 2279.8667
$ws.Range("M107").Value = -359.8667
$ws.Range("H135").Value = 851.44446
$ws.Range("I135").Value = 851.44446
$ws.Range("K135").Value = 7663.00014
$ws.Range("M135").Value = -5128.00014
$ws.Range("H138").Value = 3116.69
$ws.Range("I138").Value = 1570.6571
$ws.Range("J138").Value = 3949.1692
$ws.Range("K138").Value = 4711.971299999999
$ws.Range("L138").Value = 11847.5076
$ws.Range("M138").Value = 428.0287000000008
$ws.Range("N138").Value = -22127.5076

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 99.333336
$ws.Range("I5").Value = 99.333336
$ws.Range("K5").Value = 99.333336
$ws.Range("M5").Value = 12.666664
$ws.Range("H32").Value = 17973728
$ws.Range("I32").Value = 18141228
$ws.Range("J32").Value = 15628727
$ws.Range("K32").Value = 18141228
$ws.Range("L32").Value = 15628727
$ws.Range("M32").Value = -18140941
$ws.Range("N32").Value = -15629301
$ws.Range("H43").Value = 40780.332
$ws.Range("J43").Value = 50000
$ws.Range("L43").Value = 50000
$ws.Range("N43").Value = -50626
$ws.Range("H102").Value = 1512.1052
$ws.Range("I102").Value = 1095.0714
$ws.Range("J102").Value = 2679.8
$ws.Range("K102").Value = 1095.0714
$ws.Range("L102").Value = 2679.8
$ws.Range("M102").Value = 526.9286
$ws.Range("N102").Value = -5923.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 99.333336
$ws.Range("I4").Value = 99.333336
$ws.Range("K4").Value = 99.333336
$ws.Range("M4").Value = 15.666664
$ws.Range("H5").Value = 244.5
$ws.Range("I5").Value = 196.66667
$ws.Range("K5").Value = 196.66667
$ws.Range("M5").Value = -83.66667000000001
$ws.Range("H102").Value = 64197.6
$ws.Range("I102").Value = 10000
$ws.Range("J102").Value = 100329.336
$ws.Range("K102").Value = 10000
$ws.Range("L102").Value = 100329.336
$ws.Range("M102").Value = -6755
$ws.Range("N102").Value = -106819.336
$ws.Range("H107").Value = 1080.3846
$ws.Range("I107").Value = 900.0476
$ws.Range("K107").Value = 900.0476
$ws.Range("M107").Value = 1019.9524
$ws.Range("H134").Value = 2085450.4
$ws.Range("I134").Value = 2565663
$ws.Range("J134").Value = 4529.5
$ws.Range("K134").Value = 7696989
$ws.Range("L134").Value = 13588.5
$ws.Range("M134").Value = -7694454
$ws.Range("N134").Value = -18658.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 10843.741
$ws.Range("I86").Value = 10269.889
$ws.Range("K86").Value = 10269.889
$ws.Range("M86").Value = -9146.888999999999
$ws.Range("H89").Value = 10843.741
$ws.Range("I89").Value = 10269.889
$ws.Range("K89").Value = 51349.44499999999
$ws.Range("M89").Value = -45733.44499999999
$ws.Range("H105").Value = 1614.0952
$ws.Range("I105").Value = 1599.7333
$ws.Range("K105").Value = 1599.7333
$ws.Range("M105").Value = 147.2666999999999
$ws.Range("H130").Value = 45665.332
$ws.Range("J130").Value = 45665.332
$ws.Range("L130").Value = 45665.332
$ws.Range("N130").Value = -55705.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 111143830
$ws.Range("I4").Value = 173936080
$ws.Range("J4").Value = 66665990
$ws.Range("K4").Value = 521808240
$ws.Range("L4").Value = 199997970
$ws.Range("M4").Value = -521808128
$ws.Range("N4").Value = -199998194
$ws.Range("H15").Value = 98.75
$ws.Range("I15").Value = 22.5
$ws.Range("J15").Value = 175
$ws.Range("K15").Value = 67.5
$ws.Range("L15").Value = 525
$ws.Range("M15").Value = 72.5
$ws.Range("N15").Value = -805
$ws.Range("H82").Value = 1100
$ws.Range("I82").Value = 1100
$ws.Range("K82").Value = 3300
$ws.Range("M82").Value = -2894
$ws.Range("H85").Value = 1100
$ws.Range("I85").Value = 1100
$ws.Range("K85").Value = 3300
$ws.Range("M85").Value = -1896
$ws.Range("H86").Value = 200
$ws.Range("I86").Value = 200
$ws.Range("K86").Value = 600
$ws.Range("M86").Value = 586
$ws.Range("H89").Value = 200
$ws.Range("I89").Value = 200
$ws.Range("K89").Value = 1800
$ws.Range("M89").Value = 4128
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("M104").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 35057
$ws.Range("J52").Value = 32497.5
$ws.Range("L52").Value = 32497.5
$ws.Range("N52").Value = -33015.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1279
$ws.Range("J22").Value = 1981.4445
$ws.Range("L22").Value = 1981.4445
$ws.Range("N22").Value = -2571.4445
$ws.Range("H27").Value = 1279
$ws.Range("J27").Value = 1981.4445
$ws.Range("L27").Value = 1981.4445
$ws.Range("N27").Value = -2195.4445
$ws.Range("H45").Value = 30041
$ws.Range("I45").Value = 30041
$ws.Range("K45").Value = 30041
$ws.Range("M45").Value = -29634
$ws.Range("H46").Value = 6874.615
$ws.Range("I46").Value = 1288.5
$ws.Range("K46").Value = 1288.5
$ws.Range("M46").Value = -1100.5
$ws.Range("H82").Value = 3019.2727
$ws.Range("J82").Value = 3041.6
$ws.Range("L82").Value = 3041.6
$ws.Range("N82").Value = -3763.6
$ws.Range("H85").Value = 3019.2727
$ws.Range("J85").Value = 3041.6
$ws.Range("L85").Value = 3041.6
$ws.Range("N85").Value = -5537.6
$ws.Range("H93").Value = 50001544
$ws.Range("I93").Value = 83334750
$ws.Range("K93").Value = 83334750
$ws.Range("M93").Value = -83333502
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H133").Value = 129999
$ws.Range("J133").Value = 129999
$ws.Range("L133").Value = 129999
$ws.Range("N133").Value = -135059

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4533.56
$ws.Range("I81").Value = 917.0714
$ws.Range("K81").Value = 1834.1428
$ws.Range("M81").Value = -773.1428000000001
$ws.Range("H84").Value = 4533.56
$ws.Range("I84").Value = 917.0714
$ws.Range("K84").Value = 9170.714
$ws.Range("M84").Value = -3866.714
$ws.Range("H107").Value = 445
$ws.Range("I107").Value = 459.6875
$ws.Range("K107").Value = 1379.0625
$ws.Range("M107").Value = 540.9375
$ws.Range("H132").Value = 2238.64
$ws.Range("I132").Value = 2189.0952
$ws.Range("J132").Value = 2498.75
$ws.Range("K132").Value = 6567.285600000001
$ws.Range("L132").Value = 7496.25
$ws.Range("M132").Value = -4037.285600000001
$ws.Range("N132").Value = -12556.25
$ws.Range("H136").Value = 75865.28999999999
$ws.Range("I136").Value = 3900
$ws.Range("K136").Value = 11700
$ws.Range("M136").Value = -9150
